# Update the programme-umbrella mapping table (Table1, A1:G28) on Sheet1.
# MAVIT is split into MAVIT-RHT / MAVIT-SHA, MASY is renamed to SHA Spesialsykepleie,
# and the rows are re-grouped/re-ordered accordingly (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 'MAVIT-RHT'
$ws.Cells.Item(2,2).Value = 'MABU'
$ws.Cells.Item(2,3).Value = 'RHT Helsevitenskap'
$ws.Cells.Item(2,4).Value = 'Masterstudium i helsevitenskap - spesialisering i fysioterapi for barn og unge (D)'
$ws.Cells.Item(2,5).Value = 'M2'
$ws.Cells.Item(2,6).Value = 120
$ws.Cells.Item(2,7).Value = 0.67

$ws.Cells.Item(3,1).Value = 'MAVIT-RHT'
$ws.Cells.Item(3,2).Value = 'MAELD'
$ws.Cells.Item(3,3).Value = 'RHT Helsevitenskap'
$ws.Cells.Item(3,4).Value = 'Masterstudium i helsevitenskap - spesialisering i fysioterapi for eldre personer (D)'
$ws.Cells.Item(3,5).Value = 'M2'
$ws.Cells.Item(3,6).Value = 120
$ws.Cells.Item(3,7).Value = 0.67

$ws.Cells.Item(4,1).Value = 'MAVIT-RHT'
$ws.Cells.Item(4,2).Value = 'MAHAB'
$ws.Cells.Item(4,3).Value = 'RHT Helsevitenskap'
$ws.Cells.Item(4,4).Value = 'Masterstudium i helsevitenskap - spesialisering i rehabilitering og habilitering (D)'
$ws.Cells.Item(4,5).Value = 'M2'
$ws.Cells.Item(4,6).Value = 120
$ws.Cells.Item(4,7).Value = 0.67

$ws.Cells.Item(5,1).Value = 'MAVIT-RHT'
$ws.Cells.Item(5,2).Value = 'MERG'
$ws.Cells.Item(5,3).Value = 'RHT Helsevitenskap'
$ws.Cells.Item(5,4).Value = 'Masterstudium i helsevitenskap - spesialisering i ergoterapi (D)'
$ws.Cells.Item(5,5).Value = 'M2'
$ws.Cells.Item(5,6).Value = 120
$ws.Cells.Item(5,7).Value = 0.67

$ws.Cells.Item(6,1).Value = 'MAVIT-RHT'
$ws.Cells.Item(6,2).Value = 'MAMUS'
$ws.Cells.Item(6,3).Value = 'RHT Helsevitenskap'
$ws.Cells.Item(6,4).Value = 'Masterstudium i helsevitenskap - spesialisering i fysioterapi for muskelskjeletthelse (D)'
$ws.Cells.Item(6,5).Value = 'M2'
$ws.Cells.Item(6,6).Value = 120
$ws.Cells.Item(6,7).Value = 0.67

$ws.Cells.Item(7,1).Value = 'MAVIT-RHT'
$ws.Cells.Item(7,2).Value = 'MAPFY'
$ws.Cells.Item(7,3).Value = 'RHT Helsevitenskap'
$ws.Cells.Item(7,4).Value = 'Masterstudium i helsevitenskap - spesialisering i psykomotorisk fysioterapi (D)'
$ws.Cells.Item(7,5).Value = 'M2'
$ws.Cells.Item(7,6).Value = 120
$ws.Cells.Item(7,7).Value = 0.67

$ws.Cells.Item(8,1).Value = 'MAVIT-SHA'
$ws.Cells.Item(8,2).Value = 'MAEMP'
$ws.Cells.Item(8,3).Value = 'SHA Helsevitenskap'
$ws.Cells.Item(8,4).Value = 'Masterstudium i helsevitenskap - spesialisering i empowerment og helsefremmende arbeid'
$ws.Cells.Item(8,5).Value = 'M2'
$ws.Cells.Item(8,6).Value = 120
$ws.Cells.Item(8,7).Value = 1

$ws.Cells.Item(9,1).Value = 'MAVIT-SHA'
$ws.Cells.Item(9,2).Value = 'MAERN'
$ws.Cells.Item(9,3).Value = 'SHA Helsevitenskap'
$ws.Cells.Item(9,4).Value = 'Masterstudium i helsevitenskap - spesialisering i ernæringskompetanse for helsepersonell'
$ws.Cells.Item(9,5).Value = 'M2'
$ws.Cells.Item(9,6).Value = 120
$ws.Cells.Item(9,7).Value = 1

$ws.Cells.Item(10,1).Value = 'MAVIT-SHA'
$ws.Cells.Item(10,2).Value = 'MAHEL'
$ws.Cells.Item(10,3).Value = 'SHA Helsevitenskap'
$ws.Cells.Item(10,4).Value = 'Masterstudium i helsevitenskap - spesialisering i helsesykepleie'
$ws.Cells.Item(10,5).Value = 'M2'
$ws.Cells.Item(10,6).Value = 120
$ws.Cells.Item(10,7).Value = 1

$ws.Cells.Item(11,1).Value = 'MAVIT-SHA'
$ws.Cells.Item(11,2).Value = 'MAKLI'
$ws.Cells.Item(11,3).Value = 'SHA Helsevitenskap'
$ws.Cells.Item(11,4).Value = 'Masterstudium i helsevitenskap - spesialisering i sykepleie - klinisk forskning og fagutvikling (D)'
$ws.Cells.Item(11,5).Value = 'M2'
$ws.Cells.Item(11,6).Value = 120
$ws.Cells.Item(11,7).Value = 1

$ws.Cells.Item(12,1).Value = 'MAVIT-SHA'
$ws.Cells.Item(12,2).Value = 'MAKRE'
$ws.Cells.Item(12,3).Value = 'SHA Helsevitenskap'
$ws.Cells.Item(12,4).Value = 'Masterstudium i helsevitenskap - spesialisering i kreftsykepleie'
$ws.Cells.Item(12,5).Value = 'M2'
$ws.Cells.Item(12,6).Value = 120
$ws.Cells.Item(12,7).Value = 1

$ws.Cells.Item(13,1).Value = 'MAVIT-SHA'
$ws.Cells.Item(13,2).Value = 'MAPHN'
$ws.Cells.Item(13,3).Value = 'SHA Helsevitenskap'
$ws.Cells.Item(13,4).Value = 'Masterstudium i helsevitenskap - spesialisering i samfunnsernæring'
$ws.Cells.Item(13,5).Value = 'M2'
$ws.Cells.Item(13,6).Value = 120
$ws.Cells.Item(13,7).Value = 1

$ws.Cells.Item(14,1).Value = 'MAVIT-SHA'
$ws.Cells.Item(14,2).Value = 'MAPSY'
$ws.Cells.Item(14,3).Value = 'SHA Helsevitenskap'
$ws.Cells.Item(14,4).Value = 'Masterstudium i helsevitenskap   tverrfaglig spesialisering i psykisk helse-, rus- og avhengighetsarbeid (D)'
$ws.Cells.Item(14,5).Value = 'M2'
$ws.Cells.Item(14,6).Value = 120
$ws.Cells.Item(14,7).Value = 0.5

$ws.Cells.Item(15,1).Value = 'MASY'
$ws.Cells.Item(15,2).Value = 'MASAK'
$ws.Cells.Item(15,3).Value = 'SHA Spesialsykepleie'
$ws.Cells.Item(15,4).Value = 'Masterstudium i spesialsykepleie til akutt og kritisk syke pasienter med spesialisering i akuttsykepleie (D)'
$ws.Cells.Item(15,5).Value = 'M2'
$ws.Cells.Item(15,6).Value = 120
$ws.Cells.Item(15,7).Value = 0.8

$ws.Cells.Item(16,1).Value = 'MASY'
$ws.Cells.Item(16,2).Value = 'MASAN'
$ws.Cells.Item(16,3).Value = 'SHA Spesialsykepleie'
$ws.Cells.Item(16,4).Value = 'Masterstudium i spesialsykepleie til akutt og kritisk syke pasienter med spesialisering i anestesisykepleie (D)'
$ws.Cells.Item(16,5).Value = 'M2'
$ws.Cells.Item(16,6).Value = 120
$ws.Cells.Item(16,7).Value = 0.8

$ws.Cells.Item(17,1).Value = 'MASY'
$ws.Cells.Item(17,2).Value = 'MASBA'
$ws.Cells.Item(17,3).Value = 'SHA Spesialsykepleie'
$ws.Cells.Item(17,4).Value = 'Masterstudium i spesialsykepleie til akutt og kritisk syke pasienter med spesialisering i barnesykepleie (D)'
$ws.Cells.Item(17,5).Value = 'M2'
$ws.Cells.Item(17,6).Value = 120
$ws.Cells.Item(17,7).Value = 0.8

$ws.Cells.Item(18,1).Value = 'MASY'
$ws.Cells.Item(18,2).Value = 'MASIN'
$ws.Cells.Item(18,3).Value = 'SHA Spesialsykepleie'
$ws.Cells.Item(18,4).Value = 'Masterstudium i spesialsykepleie til akutt og kritisk syke pasienter med spesialisering i intensivsykepleie (D)'
$ws.Cells.Item(18,5).Value = 'M2'
$ws.Cells.Item(18,6).Value = 120
$ws.Cells.Item(18,7).Value = 0.8

$ws.Cells.Item(19,1).Value = 'MASY'
$ws.Cells.Item(19,2).Value = 'MASOP'
$ws.Cells.Item(19,3).Value = 'SHA Spesialsykepleie'
$ws.Cells.Item(19,4).Value = 'Masterstudium i spesialsykepleie til akutt og kritisk syke pasienter med spesialisering i operasjonssykepleie (D)'
$ws.Cells.Item(19,5).Value = 'M2'
$ws.Cells.Item(19,6).Value = 120
$ws.Cells.Item(19,7).Value = 0.8

$ws.Cells.Item(20,1).Value = 'MASY'
$ws.Cells.Item(20,2).Value = 'MASKA'
$ws.Cells.Item(20,3).Value = 'SHA Spesialsykepleie'
$ws.Cells.Item(20,4).Value = 'Masterstudium i spesialsykepleie til akutt og kritisk syke pasienter med spesialisering i kardiologisk sykepleie (D)'
$ws.Cells.Item(20,5).Value = 'M2'
$ws.Cells.Item(20,6).Value = 120
$ws.Cells.Item(20,7).Value = 0.8

$ws.Cells.Item(21,1).Value = 'MSFH'
$ws.Cells.Item(21,2).Value = 'MISWH'
$ws.Cells.Item(21,3).Value = 'MSFH Sosialfag (H)'
$ws.Cells.Item(21,4).Value = 'Master Programme in Applied Social Sciences - Study Option International Social Welfare and Health Policy'
$ws.Cells.Item(21,5).Value = 'M2'
$ws.Cells.Item(21,6).Value = 120
$ws.Cells.Item(21,7).Value = 1

$ws.Cells.Item(22,1).Value = 'MSFH'
$ws.Cells.Item(22,2).Value = 'MSBAH'
$ws.Cells.Item(22,3).Value = 'MSFH Sosialfag (H)'
$ws.Cells.Item(22,4).Value = 'Masterstudium i sosialfag - studieretning barnevernsarbeid'
$ws.Cells.Item(22,5).Value = 'M2'
$ws.Cells.Item(22,6).Value = 120
$ws.Cells.Item(22,7).Value = 1

$ws.Cells.Item(23,1).Value = 'MSFH'
$ws.Cells.Item(23,2).Value = 'MSBH'
$ws.Cells.Item(23,3).Value = 'MSFH Sosialfag (H)'
$ws.Cells.Item(23,4).Value = 'Masterstudium i sosialfag - studieretning barnevern'
$ws.Cells.Item(23,5).Value = 'M2'
$ws.Cells.Item(23,6).Value = 120
$ws.Cells.Item(23,7).Value = 1

$ws.Cells.Item(24,1).Value = 'MSFH'
$ws.Cells.Item(24,2).Value = 'MSSH'
$ws.Cells.Item(24,3).Value = 'MSFH Sosialfag (H)'
$ws.Cells.Item(24,4).Value = 'Masterstudium i sosialfag - studieretning sosialt arbeid'
$ws.Cells.Item(24,5).Value = 'M2'
$ws.Cells.Item(24,6).Value = 120
$ws.Cells.Item(24,7).Value = 1

$ws.Cells.Item(25,1).Value = 'MSFD'
$ws.Cells.Item(25,2).Value = 'MSBAD'
$ws.Cells.Item(25,3).Value = 'MSFD Sosialfag (D)'
$ws.Cells.Item(25,4).Value = 'Masterstudium i sosialfag - studieretning barnevernsarbeid, deltid (D)'
$ws.Cells.Item(25,5).Value = 'M2'
$ws.Cells.Item(25,6).Value = 120
$ws.Cells.Item(25,7).Value = 0.67

$ws.Cells.Item(26,1).Value = 'MSFD'
$ws.Cells.Item(26,2).Value = 'MSBD'
$ws.Cells.Item(26,3).Value = 'MSFD Sosialfag (D)'
$ws.Cells.Item(26,4).Value = 'Masterstudium i sosialfag - studieretning barnevern, deltid (D)'
$ws.Cells.Item(26,5).Value = 'M2'
$ws.Cells.Item(26,6).Value = 120
$ws.Cells.Item(26,7).Value = 0.67

$ws.Cells.Item(27,1).Value = 'MSFD'
$ws.Cells.Item(27,2).Value = 'MSFBD'
$ws.Cells.Item(27,3).Value = 'MSFD Sosialfag (D)'
$ws.Cells.Item(27,4).Value = 'Masterstudium i sosialfag - studieretning familiebehandling, deltid (D)'
$ws.Cells.Item(27,5).Value = 'M2'
$ws.Cells.Item(27,6).Value = 120
$ws.Cells.Item(27,7).Value = 0.67

$ws.Cells.Item(28,1).Value = 'MSFD'
$ws.Cells.Item(28,2).Value = 'MSSD'
$ws.Cells.Item(28,3).Value = 'MSFD Sosialfag (D)'
$ws.Cells.Item(28,4).Value = 'Masterstudium i sosialfag - studieretning sosialt arbeid, deltid (D)'
$ws.Cells.Item(28,5).Value = 'M2'
$ws.Cells.Item(28,6).Value = 120
$ws.Cells.Item(28,7).Value = 0.67

# Update the active selection to match the authors final cursor position.
$ws.Range("D11").Select()
